$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('B142').Value = 48654
$ws.Range('C142').Value = 'CHO-Medimix Sandal with Eladi oils for glowing skin and natural protection Soap-125gms'
$ws.Range('E142').Value = 38.26
$ws.Range('F142').Value = -1
$ws.Range('G142').Value = -32.02

$ws.Range('B143').Value = 63902
$ws.Range('C143').Value = 'CHO-Medimix Sandal with Eladi oils for glowing skin and natural protection Soap-125gms'
$ws.Range('E143').Value = 34.04
$ws.Range('F143').Value = 2
$ws.Range('G143').Value = 64.04000000000001

$ws.Range('B154').Value = 53925
$ws.Range('C154').Value = 'COL-Colgate Zigzag Charcoal Pack of 4 Toothbrush'
$ws.Range('E154').Value = 79.37
$ws.Range('F154').Value = 1
$ws.Range('G154').Value = 66.44

$ws.Range('B155').Value = 64350
$ws.Range('C155').Value = 'COL-Colgate Zigzag Charcoal Pack of 4 Toothbrush'
$ws.Range('E155').Value = 70.63
$ws.Range('F155').Value = 101
$ws.Range('G155').Value = 6710.44

$ws.Range('B156').Value = 57756
$ws.Range('C156').Value = 'COL-Colgate Zigzag Charcoal Pack of 4 Toothbrush'
$ws.Range('E156').Value = 79.37
$ws.Range('F156').Value = -100
$ws.Range('G156').Value = -6644

$ws.Range('B271').Value = 48706
$ws.Range('C271').Value = 'HIM-GENTLE BABY SOAP 75G'
$ws.Range('E271').Value = 39.8
$ws.Range('F271').Value = -144
$ws.Range('G271').Value = -4795.2

$ws.Range('B272').Value = 64973
$ws.Range('C272').Value = 'HIM-GENTLE BABY SOAP 75G'
$ws.Range('E272').Value = 35.4
$ws.Range('F272').Value = 150
$ws.Range('G272').Value = 4995

$ws.Range('B309').Value = 61610
$ws.Range('C309').Value = 'HUL-Bru Inst Poly 50g'
$ws.Range('E309').Value = 122.71
$ws.Range('F309').Value = -58
$ws.Range('G309').Value = -5957.18

$ws.Range('B310').Value = 63565
$ws.Range('C310').Value = 'HUL-Bru Inst Poly 50g'
$ws.Range('E310').Value = 109.19
$ws.Range('F310').Value = 60
$ws.Range('G310').Value = 6162.6

$ws.Range('B342').Value = 63531
$ws.Range('C342').Value = 'HUL-Kissan Pineapple Jam 500G'
$ws.Range('E342').Value = 152.53
$ws.Range('F342').Value = 80
$ws.Range('G342').Value = 11478.4

$ws.Range('B343').Value = 57802
$ws.Range('C343').Value = 'HUL-Kissan Pineapple Jam 500G'
$ws.Range('E343').Value = 162.71
$ws.Range('F343').Value = -79
$ws.Range('G343').Value = -11334.92

$ws.Range('B344').Value = 63571
$ws.Range('C344').Value = 'HUL-Kissan Pineapple Jam 500G'
$ws.Range('E344').Value = 152.53
$ws.Range('F344').Value = 29
$ws.Range('G344').Value = 4160.92

$ws.Range('B367').Value = 63563
$ws.Range('C367').Value = 'HUL-lux advanced eventoned glow 4x100'
$ws.Range('E367').Value = 119.04
$ws.Range('F367').Value = 15
$ws.Range('G367').Value = 1679.4

$ws.Range('B368').Value = 61605
$ws.Range('C368').Value = 'HUL-lux advanced eventoned glow 4x100'
$ws.Range('E368').Value = 133.78
$ws.Range('F368').Value = -13
$ws.Range('G368').Value = -1455.48

$ws.Range('B371').Value = 61608
$ws.Range('C371').Value = 'HUL-Lux Radiant Glow 3*150g'
$ws.Range('E371').Value = 154.12
$ws.Range('F371').Value = -56
$ws.Range('G371').Value = -7224.56

$ws.Range('B372').Value = 63564
$ws.Range('C372').Value = 'HUL-Lux Radiant Glow 3*150g'
$ws.Range('E372').Value = 137.16
$ws.Range('F372').Value = 57
$ws.Range('G372').Value = 7353.57

$ws.Range('B374').Value = 60325
$ws.Range('C374').Value = 'Hul-pears pure and gentle 3x125 gm'
$ws.Range('E374').Value = 151.57
$ws.Range('F374').Value = -102
$ws.Range('G374').Value = -12939.72

$ws.Range('B375').Value = 63560
$ws.Range('C375').Value = 'Hul-pears pure and gentle 3x125 gm'
$ws.Range('E375').Value = 134.87
$ws.Range('F375').Value = 104
$ws.Range('G375').Value = 13193.44

$ws.Range('B381').Value = 62865
$ws.Range('C381').Value = 'HUL-Rap Refresh Bolt 1Kg'
$ws.Range('E381').Value = 95.34999999999999
$ws.Range('F381').Value = 151
$ws.Range('G381').Value = 12051.31

$ws.Range('B382').Value = 57817
$ws.Range('C382').Value = 'HUL-Rap Refresh Bolt 1Kg'
$ws.Range('E382').Value = 95.34999999999999
$ws.Range('F382').Value = 3
$ws.Range('G382').Value = 239.43

$ws.Range('B413').Value = 57857
$ws.Range('C413').Value = 'HUL-Surf Exl Mtc Liq Tl 1 Ltr Cp'
$ws.Range('E413').Value = 180.62
$ws.Range('F413').Value = 3
$ws.Range('G413').Value = 453.51

$ws.Range('B414').Value = 63008
$ws.Range('C414').Value = 'HUL-Surf Exl Mtc Liq Tl 1 Ltr Cp'
$ws.Range('E414').Value = 180.62
$ws.Range('F414').Value = 504
$ws.Range('G414').Value = 76189.67999999999

$ws.Range('B423').Value = 63102
$ws.Range('C423').Value = 'HUL-Vim Bar Multipack Fw 4X200G'
$ws.Range('E423').Value = 71.05
$ws.Range('F423').Value = 36
$ws.Range('G423').Value = 2140.92

$ws.Range('B424').Value = 53082
$ws.Range('C424').Value = 'HUL-VIM BAR MULTIPACK FW 4X200G'
$ws.Range('E424').Value = 71.05
$ws.Range('F424').Value = 1
$ws.Range('G424').Value = 59.47

$ws.Range('B571').Value = 65069
$ws.Range('C571').Value = 'CRE-Bourbon 100gm'
$ws.Range('E571').Value = 14.3
$ws.Range('F571').Value = 172
$ws.Range('G571').Value = 2313.4

$ws.Range('B572').Value = 53757
$ws.Range('C572').Value = 'CRE-Bourbon 100gm'
$ws.Range('E572').Value = 16.08
$ws.Range('F572').Value = -159
$ws.Range('G572').Value = -2138.55

$ws.Range('B679').Value = 53319
$ws.Range('C679').Value = 'PRI-B-50 VIMAL Copper Glass 300ML (2pc Set)'
$ws.Range('E679').Value = 310.64
$ws.Range('F679').Value = -6
$ws.Range('G679').Value = -1643.52

$ws.Range('B680').Value = 64810
$ws.Range('C680').Value = 'PRI-B-50 VIMAL Copper Glass 300ML (2pc Set)'
$ws.Range('E680').Value = 291.22
$ws.Range('F680').Value = 7
$ws.Range('G680').Value = 1917.44

$ws.Range('B701').Value = 64833
$ws.Range('C701').Value = 'Rasna 32 Glass Shikanji Nimbupani'
$ws.Range('E701').Value = 34.9
$ws.Range('F701').Value = 99
$ws.Range('G701').Value = 3250.17

$ws.Range('B702').Value = 60025
$ws.Range('C702').Value = 'Rasna 32 Glass Shikanji Nimbupani'
$ws.Range('E702').Value = 37.22
$ws.Range('F702').Value = -98
$ws.Range('G702').Value = -3217.34

$ws.Range('B712').Value = 64830
$ws.Range('C712').Value = 'Rasna Nagpur Orange (32 Glass)'
$ws.Range('E712').Value = 34.9
$ws.Range('F712').Value = 117
$ws.Range('G712').Value = 3841.11

$ws.Range('B713').Value = 60022
$ws.Range('C713').Value = 'Rasna Nagpur Orange (32 Glass)'
$ws.Range('E713').Value = 37.22
$ws.Range('F713').Value = -113
$ws.Range('G713').Value = -3709.79

$ws.Range('B864').Value = 65079
$ws.Range('C864').Value = 'Shankys Tip Top Hing Jeera Peanut/ Salted Peanut 200 Gm'
$ws.Range('E864').Value = 43.44
$ws.Range('F864').Value = 21
$ws.Range('G864').Value = 858.27

$ws.Range('B865').Value = 54751
$ws.Range('C865').Value = 'Shankys Tip Top Hing Jeera Peanut/ Salted Peanut 200 Gm'
$ws.Range('E865').Value = 46.34
$ws.Range('F865').Value = -19
$ws.Range('G865').Value = -776.53

